$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Pakistan row (row 22) ---
$ws.Range("B22").Value = 56349
$ws.Range("C22").Value = 1748
$ws.Range("D22").Value = 17482
$ws.Range("E22").Value = 37700
$ws.Range("G22").Value = 34
$ws.Range("H22").Value = 1167

# --- Update Hungria row (row 74) ---
$ws.Range("B74").Value = 3756
$ws.Range("C74").Value = 15
$ws.Range("D74").Value = 1711
$ws.Range("E74").Value = 1554
$ws.Range("G74").Value = 5
$ws.Range("H74").Value = 491

# --- Update Tailandia row (row 79) ---
$ws.Range("B79").Value = 3042
$ws.Range("C79").Value = 2
$ws.Range("D79").Value = 2928
$ws.Range("E79").Value = 57
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 57

# --- El Salvador moves up in the ranking (new case totals put it ---
# --- above Macedonia/Cuba/Gabon). Re-write rows 88-91 so that El  ---
# --- Salvador occupies row 88 with its fresh totals, and          ---
# --- Macedonia/Cuba/Gabon shift down one row each, keeping their  ---
# --- own totals unchanged.                                        ---
$ws.Range("A88").Value = "El Salvador"
$ws.Range("B88").Value = 1983
$ws.Range("C88").Value = 68
$ws.Range("D88").Value = 698
$ws.Range("E88").Value = 1250
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 35

$ws.Range("A89").Value = "Republica de Macedonia"
$ws.Range("B89").Value = 1978
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 1422
$ws.Range("E89").Value = 443
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 113

$ws.Range("A90").Value = "Cuba"
$ws.Range("B90").Value = 1941
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 1689
$ws.Range("E90").Value = 170
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 82

$ws.Range("A91").Value = "Gabon"
$ws.Range("B91").Value = 1934
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 459
$ws.Range("E91").Value = 1463
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 12

# --- Update the "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 07:05"
